# Sequence diagram touch-up:
#  - rename deletePerson(p) -> deleteTask(p) call-out text
#  - shrink several call-out text boxes from 14pt to 12pt and tighten
#    their auto-fit box heights to match (215444 EMU -> 184666 EMU,
#    i.e. 16.96409pt -> 14.54063pt)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU-per-point constant used to convert the target <a:ext cy="..."/> values
# (taken from the canonical OOXML) into the point units the Shape.Height
# property expects.
$EMU_PER_PT = 12700
$newHeightPt = 184666 / $EMU_PER_PT

# --- Shape "TextBox 25" (id 26) : execute("delete 1") ----------------------
$sh = $s.Shapes.Item(13)
$sh.TextFrame.TextRange.Font.Size = 12
$sh.Height = $newHeightPt

# --- Shape "TextBox 28" (id 29) : execute() ---------------------------------
$sh = $s.Shapes.Item(15)
$sh.TextFrame.TextRange.Font.Size = 12
$sh.Height = $newHeightPt

# --- Shape "TextBox 77" (id 78) : deletePerson(p) -> deleteTask(p) ---------
$sh = $s.Shapes.Item(27)
$tr = $sh.TextFrame.TextRange
# Rename the first run's text ("deletePerson" -> "deleteTask").
$tr.Characters(1, 12).Text = "deleteTask"
# Resize every run (and implicitly split "(p)" into "(p" + ")" once we touch
# the trailing character on its own).
$tr.Font.Size = 12
$tr.Characters($tr.Length, 1).Text = ")"
$sh.Height = $newHeightPt

# --- Shape "TextBox 79" (id 80) : parse("delete 1") -------------------------
$sh = $s.Shapes.Item(29)
$sh.TextFrame.TextRange.Font.Size = 12
$sh.Height = $newHeightPt

# --- Shape "TextBox 81" (id 82) : result ------------------------------------
$sh = $s.Shapes.Item(30)
$sh.TextFrame.TextRange.Font.Size = 12
$sh.Height = $newHeightPt

# --- Shape "TextBox 82" (id 83) : result ------------------------------------
$sh = $s.Shapes.Item(31)
$sh.TextFrame.TextRange.Font.Size = 12
$sh.Height = $newHeightPt
